$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set text cells via NumberFormat "@" to prevent Excel from auto-parsing
# numeric-looking strings (e.g. "301.24", "0.999") into floating point
# numbers, which would lose the exact textual representation (trailing
# zeros, thousand-dot grouping, etc.) used by this sheet. After setting
# the value we reset the style back to "Normal" so no stray number format
# is left behind on the cell (matching the original plain inline-string cells).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "42.834.55"
$ws.Range("E2").Value = "  -1.52%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.301.10"
$ws.Range("E3").Value = "  -0.49%  "

# Row 4
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
Set-TextValue $ws.Range("D5") "301.24"
$ws.Range("E5").Value = "  -2.43%  "

# Row 6
Set-TextValue $ws.Range("D6") "98.40"
$ws.Range("E6").Value = "  -7.11%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.503"
$ws.Range("E7").Value = "  -4.33%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.500"
$ws.Range("E9").Value = "  -3.70%  "

# Row 10
Set-TextValue $ws.Range("D10") "34.59"
$ws.Range("E10").Value = "  -4.74%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.0788"
$ws.Range("E11").Value = "  -3.10%  "

# Row 12
$ws.Range("E12").Value = "  +0.25%  "

# Row 13
Set-TextValue $ws.Range("B13") "Polkadot"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D13") "6.67"
$ws.Range("E13").Value = "  -4.38%  "

# Row 14
Set-TextValue $ws.Range("B14") "WrappedliquidstakedEther2.0"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D14") "2.652.02"
$ws.Range("E14").Value = "  -0.63%  "

# Row 15
$ws.Range("E15").Value = "  +0.36%  "

# Row 16
Set-TextValue $ws.Range("D16") "2.270.15"
$ws.Range("E16").Value = "  -2.26%  "

# Row 17
Set-TextValue $ws.Range("D17") "0.789"
$ws.Range("E17").Value = "  -1.58%  "

# Row 18
Set-TextValue $ws.Range("D18") "42.719.62"
$ws.Range("E18").Value = "  -1.65%  "

# Row 19
Set-TextValue $ws.Range("D19") "11.57"
$ws.Range("E19").Value = "  -3.03%  "

# Row 20
Set-TextValue $ws.Range("D20") "0.0₃0896"
$ws.Range("E20").Value = "  -2.86%  "

# Row 21
Set-TextValue $ws.Range("D21") "5.99"
$ws.Range("E21").Value = "  -3.34%  "

# Row 22
Set-TextValue $ws.Range("D22") "67.58"
$ws.Range("E22").Value = "  -0.74%  "

# Row 23
Set-TextValue $ws.Range("D23") "235.65"
$ws.Range("E23").Value = "  -2.34%  "

# Row 24
$ws.Range("E24").Value = "  -4.85%  "

# Row 25
$ws.Range("E25").Value = "  -4.58%  "

# Row 26
$ws.Range("E26").Value = "  +0.25%  "

# Row 27
Set-TextValue $ws.Range("D27") "24.53"
$ws.Range("E27").Value = "  -1.98%  "

# Row 28
$ws.Range("E28").Value = "  -2.77%  "

# Row 29
Set-TextValue $ws.Range("D29") "34.08"
$ws.Range("E29").Value = "  -6.83%  "

# Row 30
Set-TextValue $ws.Range("D30") "163.77"
$ws.Range("E30").Value = "  +0.25%  "

# Row 31
Set-TextValue $ws.Range("D31") "9.07"
$ws.Range("E31").Value = "  -5.44%  "

# Row 32
Set-TextValue $ws.Range("D32") "0.998"
$ws.Range("E32").Value = "  -0.17%  "

# Row 33
Set-TextValue $ws.Range("D33") "4.98"
$ws.Range("E33").Value = "  -4.94%  "

# Row 35
Set-TextValue $ws.Range("D35") "4.43"
$ws.Range("E35").Value = "  -3.97%  "

# Row 36
Set-TextValue $ws.Range("B36") "Hedera"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D36") "0.0692"
$ws.Range("E36").Value = "  -5.87%  "

# Row 37
Set-TextValue $ws.Range("B37") "Celestia"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Range("D37") "16.52"
$ws.Range("E37").Value = "  -9.41%  "

# Row 38
Set-TextValue $ws.Range("D38") "2.87"
$ws.Range("E38").Value = "  -4.89%  "

# Row 39
Set-TextValue $ws.Range("D39") "1.78"
$ws.Range("E39").Value = "  -4.45%  "

# Row 40
$ws.Range("E40").Value = "  -5.42%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.110"
$ws.Range("E41").Value = "  -4.68%  "

# Row 42
$ws.Range("E42").Value = "  +0.60%  "

# Row 43
Set-TextValue $ws.Range("D43") "1.954.04"
$ws.Range("E43").Value = "  -0.41%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.0279"
$ws.Range("E44").Value = "  -3.74%  "

# Row 45
Set-TextValue $ws.Range("D45") "18.29"
$ws.Range("E45").Value = "  -2.58%  "

# Row 46
Set-TextValue $ws.Range("D46") "10.11"
$ws.Range("E46").Value = "  -1.84%  "

# Row 47
Set-TextValue $ws.Range("D47") "2.86"
$ws.Range("E47").Value = "  -6.69%  "

# Row 48
Set-TextValue $ws.Range("D48") "54.76"
$ws.Range("E48").Value = "  -5.71%  "

# Row 49
Set-TextValue $ws.Range("D49") "2.524.61"
$ws.Range("E49").Value = "  -0.51%  "

# Row 50
$ws.Range("E50").Value = "  -5.17%  "

# Row 51
Set-TextValue $ws.Range("D51") "4.66"
$ws.Range("E51").Value = "  -1.76%  "
